$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -76.923076923076
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = -36.666666666666
$ws.Range("L15").Value = 11.764705882352
$ws.Range("M15").Value = 5.555555555555
$ws.Range("N15").Value = -67.241379310344
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 99
$ws.Range("J16").Value = 96
$ws.Range("K16").Value = 3.125
$ws.Range("L16").Value = 1.020408163265
$ws.Range("M16").Value = -58.227848101265
$ws.Range("N16").Value = -87.275064267352
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 246
$ws.Range("J17").Value = 288
$ws.Range("K17").Value = -14.583333333333
$ws.Range("L17").Value = -12.765957446808
$ws.Range("M17").Value = 2.074688796680
$ws.Range("N17").Value = -51.669941060903
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 86
$ws.Range("K18").Value = -38.372093023255
$ws.Range("L18").Value = -37.647058823529
$ws.Range("M18").Value = -79.215686274509
$ws.Range("N18").Value = -92.779291553133
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -50
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 3.571428571428
$ws.Range("I19").Value = 283
$ws.Range("J19").Value = 308
$ws.Range("K19").Value = -8.116883116883
$ws.Range("L19").Value = 23.043478260869
$ws.Range("M19").Value = -32.458233890214
$ws.Range("N19").Value = -90.576090576090
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -38.461538461538
$ws.Range("I20").Value = 143
$ws.Range("J20").Value = 172
$ws.Range("K20").Value = -16.860465116279
$ws.Range("L20").Value = 14.4
$ws.Range("M20").Value = -23.936170212766
$ws.Range("N20").Value = -87.400881057268
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -30
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -15.094339622641
$ws.Range("I21").Value = 846
$ws.Range("J21").Value = 988
$ws.Range("K21").Value = -14.372469635627
$ws.Range("L21").Value = -0.118063754427
$ws.Range("M21").Value = -38.293216630196
$ws.Range("N21").Value = -86.435786435786
$ws.Range("I23").Value = 10
$ws.Range("K23").Value = 233.333333333333
$ws.Range("L23").Value = -41.176470588235
$ws.Range("M23").Value = 25
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 20.879120879120
$ws.Range("I24").Value = 838
$ws.Range("J24").Value = 907
$ws.Range("K24").Value = -7.607497243660
$ws.Range("L24").Value = 31.142410015649
$ws.Range("M24").Value = 18.028169014084
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 46
$ws.Range("H25").Value = 23.913043478260
$ws.Range("I25").Value = 460
$ws.Range("J25").Value = 369
$ws.Range("K25").Value = 24.661246612466
$ws.Range("L25").Value = 42.414860681114
$ws.Range("M25").Value = -19.860627177700
$ws.Range("E26").Value = -50
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -85.714285714285
$ws.Range("I26").Value = 31
$ws.Range("J26").Value = 41
$ws.Range("K26").Value = -24.390243902439
$ws.Range("L26").Value = -6.060606060606
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -80
$ws.Range("J28").Value = 27
$ws.Range("K28").Value = -40.740740740740
$ws.Range("M28").Value = -61.904761904761
$ws.Range("N28").Value = -83.333333333333
$ws.Range("J29").Value = 21
$ws.Range("K29").Value = -47.619047619047
$ws.Range("M29").Value = -68.571428571428
$ws.Range("N29").Value = -86.904761904761

# --- Cells that change type (numeric <-> text) need both a value write and a
#     style/number-format fix-up, since Excel keeps the previous cell format
#     when only .Value is assigned. We copy the number format from a neighboring
#     cell that already has the desired look and paste formats only (xlPasteFormats).
$ws.Range("C15").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C23").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C23").PasteSpecial(-4122)

$ws.Range("C26").Value = 1
$ws.Range("C16").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("D28").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = -100
$ws.Range("E26").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("D29").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").Value = -100
$ws.Range("E26").Copy()
$ws.Range("E29").PasteSpecial(-4122)

